# Scheduled runner refresh: update market-board derived figures
# (currentAveragePrice/NQ/HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# for the Leve rows whose prices moved since the last fetch.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 7405.9287
$ws.Cells.Item(9, 9).Value = 321.8
$ws.Cells.Item(9, 11).Value = 321.8
$ws.Cells.Item(9, 13).Value = -152.8

$ws.Cells.Item(111, 8).Value = 1186.8667
$ws.Cells.Item(111, 9).Value = 1169
$ws.Cells.Item(111, 11).Value = 3507
$ws.Cells.Item(111, 13).Value = -440

$ws.Cells.Item(132, 8).Value = 2592.4614
$ws.Cells.Item(132, 9).Value = 2604.5715
$ws.Cells.Item(132, 11).Value = 7813.7145
$ws.Cells.Item(132, 13).Value = -5283.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1392.4286
$ws.Cells.Item(2, 9).Value = 1863.6364
$ws.Cells.Item(2, 10).Value = 1087.5294
$ws.Cells.Item(2, 11).Value = 1863.6364
$ws.Cells.Item(2, 12).Value = 1087.5294
$ws.Cells.Item(2, 13).Value = -1750.6364
$ws.Cells.Item(2, 14).Value = -1313.5294

$ws.Cells.Item(61, 8).Value = 3092.349
$ws.Cells.Item(61, 10).Value = 4343.25
$ws.Cells.Item(61, 12).Value = 4343.25
$ws.Cells.Item(61, 14).Value = -4767.25

$ws.Cells.Item(116, 8).Value = 1392.4286
$ws.Cells.Item(116, 9).Value = 1863.6364
$ws.Cells.Item(116, 10).Value = 1087.5294
$ws.Cells.Item(116, 11).Value = 1863.6364
$ws.Cells.Item(116, 12).Value = 1087.5294
$ws.Cells.Item(116, 13).Value = 430.3635999999999
$ws.Cells.Item(116, 14).Value = -5675.529399999999

$ws.Cells.Item(136, 8).Value = 3092.349
$ws.Cells.Item(136, 10).Value = 4343.25
$ws.Cells.Item(136, 12).Value = 13029.75
$ws.Cells.Item(136, 14).Value = -18129.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1392.4286
$ws.Cells.Item(3, 9).Value = 1863.6364
$ws.Cells.Item(3, 10).Value = 1087.5294
$ws.Cells.Item(3, 11).Value = 1863.6364
$ws.Cells.Item(3, 12).Value = 1087.5294
$ws.Cells.Item(3, 13).Value = -1749.6364
$ws.Cells.Item(3, 14).Value = -1315.5294

$ws.Cells.Item(20, 8).Value = 2223.9375
$ws.Cells.Item(20, 9).Value = 1728.8096
$ws.Cells.Item(20, 10).Value = 3169.182
$ws.Cells.Item(20, 11).Value = 1728.8096
$ws.Cells.Item(20, 12).Value = 3169.182
$ws.Cells.Item(20, 13).Value = -1481.8096
$ws.Cells.Item(20, 14).Value = -3663.182

$ws.Cells.Item(22, 8).Value = 323.1111
$ws.Cells.Item(22, 9).Value = 303
$ws.Cells.Item(22, 10).Value = 363.33334
$ws.Cells.Item(22, 11).Value = 303
$ws.Cells.Item(22, 12).Value = 363.33334
$ws.Cells.Item(22, 13).Value = -130
$ws.Cells.Item(22, 14).Value = -709.33334

$ws.Cells.Item(86, 8).Value = 66669380
$ws.Cells.Item(86, 9).Value = 125002056
$ws.Cells.Item(86, 10).Value = 3469.4285
$ws.Cells.Item(86, 11).Value = 125002056
$ws.Cells.Item(86, 12).Value = 3469.4285
$ws.Cells.Item(86, 13).Value = -125000933
$ws.Cells.Item(86, 14).Value = -5715.4285

$ws.Cells.Item(89, 8).Value = 66669380
$ws.Cells.Item(89, 9).Value = 125002056
$ws.Cells.Item(89, 10).Value = 3469.4285
$ws.Cells.Item(89, 11).Value = 625010280
$ws.Cells.Item(89, 12).Value = 17347.1425
$ws.Cells.Item(89, 13).Value = -625004664
$ws.Cells.Item(89, 14).Value = -28579.1425

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 52086270
$ws.Cells.Item(58, 9).Value = 41669076
$ws.Cells.Item(58, 11).Value = 41669076
$ws.Cells.Item(58, 13).Value = -41668873

$ws.Cells.Item(62, 8).Value = 294994
$ws.Cells.Item(62, 9).Value = 342495.34
$ws.Cells.Item(62, 11).Value = 342495.34
$ws.Cells.Item(62, 13).Value = -341871.34

$ws.Cells.Item(65, 8).Value = 294994
$ws.Cells.Item(65, 9).Value = 342495.34
$ws.Cells.Item(65, 11).Value = 1712476.7
$ws.Cells.Item(65, 13).Value = -1709356.7

$ws.Cells.Item(86, 8).Value = 9977.115
$ws.Cells.Item(86, 9).Value = 8005.5386
$ws.Cells.Item(86, 10).Value = 11948.692
$ws.Cells.Item(86, 11).Value = 8005.5386
$ws.Cells.Item(86, 12).Value = 11948.692
$ws.Cells.Item(86, 13).Value = -6882.5386
$ws.Cells.Item(86, 14).Value = -14194.692

$ws.Cells.Item(89, 8).Value = 9977.115
$ws.Cells.Item(89, 9).Value = 8005.5386
$ws.Cells.Item(89, 10).Value = 11948.692
$ws.Cells.Item(89, 11).Value = 40027.693
$ws.Cells.Item(89, 12).Value = 59743.45999999999
$ws.Cells.Item(89, 13).Value = -34411.693
$ws.Cells.Item(89, 14).Value = -70975.45999999999

$ws.Cells.Item(132, 8).Value = 4351.6294
$ws.Cells.Item(132, 9).Value = 3156.739
$ws.Cells.Item(132, 10).Value = 11222.25
$ws.Cells.Item(132, 11).Value = 9470.217000000001
$ws.Cells.Item(132, 12).Value = 33666.75
$ws.Cells.Item(132, 13).Value = -6940.217000000001
$ws.Cells.Item(132, 14).Value = -38726.75

$ws.Cells.Item(134, 8).Value = 3253.48
$ws.Cells.Item(134, 9).Value = 3496
$ws.Cells.Item(134, 10).Value = 2485.5
$ws.Cells.Item(134, 11).Value = 10488
$ws.Cells.Item(134, 12).Value = 7456.5
$ws.Cells.Item(134, 13).Value = -7953
$ws.Cells.Item(134, 14).Value = -12526.5

$ws.Cells.Item(136, 8).Value = 52086270
$ws.Cells.Item(136, 9).Value = 41669076
$ws.Cells.Item(136, 11).Value = 125007228
$ws.Cells.Item(136, 13).Value = -125004678

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5451.778
$ws.Cells.Item(70, 9).Value = 5430.1665
$ws.Cells.Item(70, 11).Value = 5430.1665
$ws.Cells.Item(70, 13).Value = -5160.1665

$ws.Cells.Item(73, 8).Value = 5451.778
$ws.Cells.Item(73, 9).Value = 5430.1665
$ws.Cells.Item(73, 11).Value = 5430.1665
$ws.Cells.Item(73, 13).Value = -4494.1665

$ws.Cells.Item(132, 8).Value = 11943.444
$ws.Cells.Item(132, 9).Value = 11165.5
$ws.Cells.Item(132, 11).Value = 33496.5
$ws.Cells.Item(132, 13).Value = -30966.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 16274.625
$ws.Cells.Item(7, 9).Value = 22161.8
$ws.Cells.Item(7, 10).Value = 6462.6665
$ws.Cells.Item(7, 11).Value = 22161.8
$ws.Cells.Item(7, 12).Value = 6462.6665
$ws.Cells.Item(7, 13).Value = -22049.8
$ws.Cells.Item(7, 14).Value = -6686.6665

$ws.Cells.Item(61, 8).Value = 1688.6818
$ws.Cells.Item(61, 9).Value = 1662.7222
$ws.Cells.Item(61, 11).Value = 1662.7222
$ws.Cells.Item(61, 13).Value = -1460.7222

$ws.Cells.Item(113, 8).Value = 1688.6818
$ws.Cells.Item(113, 9).Value = 1662.7222
$ws.Cells.Item(113, 11).Value = 1662.7222
$ws.Cells.Item(113, 13).Value = 507.2778000000001

$ws.Cells.Item(126, 8).Value = 16274.625
$ws.Cells.Item(126, 9).Value = 22161.8
$ws.Cells.Item(126, 10).Value = 6462.6665
$ws.Cells.Item(126, 11).Value = 66485.39999999999
$ws.Cells.Item(126, 12).Value = 19387.9995
$ws.Cells.Item(126, 13).Value = -64015.39999999999
$ws.Cells.Item(126, 14).Value = -24327.9995

$ws.Cells.Item(136, 8).Value = 27343.166
$ws.Cells.Item(136, 9).Value = 2320.0715
$ws.Cells.Item(136, 10).Value = 114924
$ws.Cells.Item(136, 11).Value = 6960.2145
$ws.Cells.Item(136, 12).Value = 344772
$ws.Cells.Item(136, 13).Value = -4410.2145
$ws.Cells.Item(136, 14).Value = -349872

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 772.5909
$ws.Cells.Item(113, 9).Value = 753.19354
$ws.Cells.Item(113, 11).Value = 2259.58062
$ws.Cells.Item(113, 13).Value = -89.58061999999973

$ws.Cells.Item(132, 8).Value = 22174
$ws.Cells.Item(132, 9).Value = 22174
$ws.Cells.Item(132, 11).Value = 66522
$ws.Cells.Item(132, 13).Value = -63992

$ws.Cells.Item(136, 8).Value = 41424860
$ws.Cells.Item(136, 9).Value = 3499849.8
$ws.Cells.Item(136, 10).Value = 250012420
$ws.Cells.Item(136, 11).Value = 10499549.4
$ws.Cells.Item(136, 12).Value = 750037260
$ws.Cells.Item(136, 13).Value = -10496999.4
$ws.Cells.Item(136, 14).Value = -750042360
